$d = $word.ActiveDocument

# --- Paragraph 1: title ---
$p1 = $d.Paragraphs(1)
$p1.Format.SpaceAfter = 24
$p1.Range.Find.Execute("Okos panzió", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Okospanzió", 2)

# --- Fix typos in the big paragraph (still paragraph 3 at this point) ---
$big = $d.Paragraphs(3).Range
$big.Find.Execute("emelett", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "mellette", 2)
$big = $d.Paragraphs(3).Range
$big.Find.Execute("mobilalkmazás", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "mobilalkalmazás", 2)
$big = $d.Paragraphs(3).Range
$big.Find.Execute("adminisztrácós", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "adminisztrációs", 2)
$big = $d.Paragraphs(3).Range
$big.Find.Execute("felhasznákó", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "felhasználó", 2)

# --- Move the corrected text into paragraph 2, then blank out paragraph 3 ---
$p3 = $d.Paragraphs(3)
# paragraph 3 text without its trailing paragraph mark
$bodyRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$bodyText = $bodyRange.Text

$p2 = $d.Paragraphs(2)
$p2.Range.Text = $bodyText
$p2.Range.Font.Size = 12

# Now locate paragraph 3 again (content shifted) and clear all but the paragraph mark,
# keep the bookmark that was there.
$p3 = $d.Paragraphs(3)
$clearRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$clearRange.Text = ""
$p3.Range.Font.Size = 12

# --- Remove the final (4th) empty "Listaszerű bekezdés" paragraph ---
$p4 = $d.Paragraphs(4)
$p4.Range.Delete()
